$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.930.84'
$ws.Range('E2').Value = '  +4.18%  '

# Row 3
$ws.Range('D3').Value = '3.249.63'
$ws.Range('E3').Value = '  +2.03%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').Value = "'578.23"
$ws.Range('E5').Value = '  +2.87%  '

# Row 6
$ws.Range('D6').Value = "'177.03"
$ws.Range('E6').Value = '  +2.68%  '

# Row 7
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = "'0.603"
$ws.Range('E8').Value = '  +0.22%  '

# Row 9
$ws.Range('D9').Value = '3.247.58'
$ws.Range('E9').Value = '  +1.92%  '

# Row 10
$ws.Range('E10').Value = '  +3.95%  '

# Row 11
$ws.Range('D11').Value = "'6.71"
$ws.Range('E11').Value = '  +1.50%  '

# Row 12
$ws.Range('D12').Value = "'0.407"
$ws.Range('E12').Value = '  +2.38%  '

# Row 13
$ws.Range('D13').Value = '3.815.40'
$ws.Range('E13').Value = '  +2.15%  '

# Row 14
$ws.Range('E14').Value = '  +1.48%  '

# Row 15
$ws.Range('D15').Value = "'27.87"
$ws.Range('E15').Value = '  +1.22%  '

# Row 16
$ws.Range('D16').Value = '66.924.10'
$ws.Range('E16').Value = '  +4.16%  '

# Row 17
$ws.Range('E17').Value = '  +2.85%  '

# Row 18
$ws.Range('D18').Value = '3.251.66'
$ws.Range('E18').Value = '  +2.36%  '

# Row 19
$ws.Range('D19').Value = "'5.79"
$ws.Range('E19').Value = '  +2.05%  '

# Row 20
$ws.Range('E20').Value = '  +1.97%  '

# Row 21
$ws.Range('D21').Value = "'368.88"
$ws.Range('E21').Value = '  +4.25%  '

# Row 22
$ws.Range('E22').Value = '  +4.41%  '

# Row 23
$ws.Range('E23').Value = '  +0.07%  '

# Row 24
$ws.Range('D24').Value = "'70.63"
$ws.Range('E24').Value = '  +1.83%  '

# Row 25
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.387.29'
$ws.Range('E25').Value = '  +2.42%  '

# Row 26
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').Value = "'0.506"
$ws.Range('E26').Value = '  +0.69%  '

# Row 27
$ws.Range('E27').Value = '  -0.43%  '

# Row 28
$ws.Range('E28').Value = '  +2.99%  '

# Row 29
$ws.Range('E29').Value = '  +1.63%  '

# Row 30
$ws.Range('E30').Value = '  -0.07%  '

# Row 31
$ws.Range('E31').Value = '  +4.61%  '

# Row 32
$ws.Range('D32').Value = "'5.63"
$ws.Range('E32').Value = '  -0.45%  '

# Row 33
$ws.Range('D33').Value = "'22.46"
$ws.Range('E33').Value = '  +1.35%  '

# Row 34
$ws.Range('E34').Value = '  -0.08%  '

# Row 35
$ws.Range('D35').Value = "'173.64"
$ws.Range('E35').Value = '  +10.52%  '

# Row 36
$ws.Range('E36').Value = '  +2.47%  '

# Row 37
$ws.Range('D37').Value = "'6.75"
$ws.Range('E37').Value = '  +1.78%  '

# Row 39
$ws.Range('D39').Value = "'0.852"
$ws.Range('E39').Value = '  +7.01%  '

# Row 40
$ws.Range('E40').Value = '  +9.26%  '

# Row 41
$ws.Range('D41').Value = "'26.79"
$ws.Range('E41').Value = '  +2.53%  '

# Row 42
$ws.Range('E42').Value = '  +1.47%  '

# Row 43
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'6.43"
$ws.Range('E43').Value = '  +6.53%  '

# Row 44
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.720.58'
$ws.Range('E44').Value = '  +2.29%  '

# Row 45
$ws.Range('E45').Value = '  +2.78%  '

# Row 46
$ws.Range('D46').Value = "'40.40"
$ws.Range('E46').Value = '  +4.00%  '

# Row 47
$ws.Range('D47').Value = "'0.0672"
$ws.Range('E47').Value = '  +2.75%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = "'24.61"
$ws.Range('E48').Value = '  +3.58%  '

# Row 49
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = "'334.63"
$ws.Range('E49').Value = '  +1.45%  '

# Row 50
$ws.Range('D50').Value = "'0.0278"
$ws.Range('E50').Value = '  +2.69%  '

# Row 51
$ws.Range('D51').Value = "'0.103"
$ws.Range('E51').Value = '  +2.18%  '
